# Auto-generated COM-interop script applying "added util hrs and total hrs display"
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("F15")
$ws2 = $wb.Worksheets.Item("F20")
$ws3 = $wb.Worksheets.Item("TRX")

# --- F15 (sheet1): rows 5-18 ---
$ws1.Range("G5").Value = 123
$ws1.Range("H5").Value = 130
$ws1.Range("G6").Value = 230
$ws1.Range("H6").Value = 231
$ws1.Range("G7").Value = 22
$ws1.Range("H7").Value = 34
$ws1.Range("G8").Value = 109
$ws1.Range("H8").Value = 119
$ws1.Range("G9").Value = 110
$ws1.Range("H9").Value = 120
$ws1.Range("G10").Value = 111
$ws1.Range("H10").Value = 121
$ws1.Range("G11").Value = 112
$ws1.Range("H11").Value = 122
$ws1.Range("G12").Value = 113
$ws1.Range("H12").Value = 123
$ws1.Range("G13").Value = 114
$ws1.Range("H13").Value = 124
$ws1.Range("G14").Value = 115
$ws1.Range("H14").Value = 125
$ws1.Range("G15").Value = 763
$ws1.Range("H15").Value = 799
$ws1.Range("G16").Value = 77
$ws1.Range("H16").Value = 80
$ws1.Range("G17").Value = 90
$ws1.Range("H17").Value = 95
$ws1.Range("G18").Value = 120
$ws1.Range("H18").Value = 125

$ws1.Range("I5").Formula = "=H5-G5"
$ws1.Range("I6:I18").Formula = "=H6-G6"

# --- F20 (sheet2): rows 5-15 have data; rows 16-18 stay blank ---
$ws2.Range("G5").Value = 123
$ws2.Range("H5").Value = 130
$ws2.Range("G6").Value = 230
$ws2.Range("H6").Value = 231
$ws2.Range("G7").Value = 22
$ws2.Range("H7").Value = 34
$ws2.Range("G8").Value = 109
$ws2.Range("H8").Value = 119
$ws2.Range("G9").Value = 110
$ws2.Range("H9").Value = 120
$ws2.Range("G10").Value = 111
$ws2.Range("H10").Value = 121
$ws2.Range("G11").Value = 112
$ws2.Range("H11").Value = 122
$ws2.Range("G12").Value = 113
$ws2.Range("H12").Value = 123
$ws2.Range("G13").Value = 114
$ws2.Range("H13").Value = 124
$ws2.Range("G14").Value = 115
$ws2.Range("H14").Value = 125
$ws2.Range("G15").Value = 763
$ws2.Range("H15").Value = 799

$ws2.Range("I5:I18").Formula = "=H5-G5"
$ws2.Range("I16:I18").ClearContents()

# --- TRX (sheet3): rows 5-15 have data; rows 16-18 stay blank ---
$ws3.Range("G5").Value = 123
$ws3.Range("H5").Value = 130
$ws3.Range("G6").Value = 230
$ws3.Range("H6").Value = 231
$ws3.Range("G7").Value = 22
$ws3.Range("H7").Value = 34
$ws3.Range("G8").Value = 109
$ws3.Range("H8").Value = 119
$ws3.Range("G9").Value = 110
$ws3.Range("H9").Value = 120
$ws3.Range("G10").Value = 111
$ws3.Range("H10").Value = 121
$ws3.Range("G11").Value = 112
$ws3.Range("H11").Value = 122
$ws3.Range("G12").Value = 113
$ws3.Range("H12").Value = 123
$ws3.Range("G13").Value = 114
$ws3.Range("H13").Value = 124
$ws3.Range("G14").Value = 115
$ws3.Range("H14").Value = 125
$ws3.Range("G15").Value = 763
$ws3.Range("H15").Value = 799

$ws3.Range("I5:I18").Formula = "=H5-G5"
$ws3.Range("I16:I18").ClearContents()

# --- Selection / active sheet bookkeeping (order matters: do last the sheet that must end up active) ---
$ws1.Range("G5:I18").Select()
$ws3.Range("G16:J22").Select()
$ws2.Range("P18").Select()

